# Auto-generated Excel COM-interop script
# Applies per-cell numeric updates (and a few cell adds/removals) to the
# Leve profit-tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 1679.4
$ws.Range("I76").Value = 1349.25
$ws.Range("K76").Value = 1349.25
$ws.Range("M76").Value = -1034.25

$ws.Range("H79").Value = 1679.4
$ws.Range("I79").Value = 1349.25
$ws.Range("K79").Value = 1349.25
$ws.Range("M79").Value = -257.25

$ws.Range("H88").Value = 741.5
$ws.Range("J88").Value = 741.5
$ws.Range("L88").Value = 741.5
$ws.Range("N88").Value = -1553.5

$ws.Range("H91").Value = 741.5
$ws.Range("J91").Value = 741.5
$ws.Range("L91").Value = 741.5
$ws.Range("N91").Value = -3549.5

$ws.Range("H98").Value = 4501
$ws.Range("J98").Value = 7000
$ws.Range("L98").Value = 7000
$ws.Range("N98").Value = -9996

$ws.Range("H122").Value = 4501
$ws.Range("J122").Value = 7000
$ws.Range("L122").Value = 21000
$ws.Range("N122").Value = -25900

$ws.Range("H132").Value = 1237.7059
$ws.Range("I132").Value = 1237.7059
$ws.Range("K132").Value = 3713.1177
$ws.Range("M132").Value = -1183.1177

$ws.Range("H135").Value = 1125.1875
$ws.Range("I135").Value = 619.0909
$ws.Range("K135").Value = 5571.8181
$ws.Range("M135").Value = -3036.8181

$ws.Range("H137").Value = 3697.8333
$ws.Range("I137").Value = 3387.4
$ws.Range("K137").Value = 10162.2
$ws.Range("M137").Value = -7612.200000000001

$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11208
$ws.Range("I32").Value = 11208
$ws.Range("K32").Value = 11208
$ws.Range("M32").Value = -10921

$ws.Range("H74").Value = 24598.643
$ws.Range("I74").Value = 23414
$ws.Range("K74").Value = 23414
$ws.Range("M74").Value = -22540

$ws.Range("H77").Value = 24598.643
$ws.Range("I77").Value = 23414
$ws.Range("K77").Value = 117070
$ws.Range("M77").Value = -112702

$ws.Range("H122").Value = 2618.4546
$ws.Range("I122").Value = 2680.3
$ws.Range("K122").Value = 8040.900000000001
$ws.Range("M122").Value = -5590.900000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5990
$ws.Range("I20").Value = 1320.1666
$ws.Range("K20").Value = 1320.1666
$ws.Range("M20").Value = -1073.1666

$ws.Range("H86").Value = 2647.25
$ws.Range("I86").Value = 2929.6667
$ws.Range("J86").Value = 1800
$ws.Range("K86").Value = 2929.6667
$ws.Range("L86").Value = 1800
$ws.Range("M86").Value = -1806.6667
$ws.Range("N86").Value = -4046

$ws.Range("H89").Value = 2647.25
$ws.Range("I89").Value = 2929.6667
$ws.Range("J89").Value = 1800
$ws.Range("K89").Value = 14648.3335
$ws.Range("L89").Value = 9000
$ws.Range("M89").Value = -9032.333500000001
$ws.Range("N89").Value = -20232

$ws.Range("H105").Value = 2777.2144
$ws.Range("I105").Value = 2777.2144
$ws.Range("K105").Value = 2777.2144
$ws.Range("M105").Value = -1030.2144

$ws.Range("H107").Value = 1207.25
$ws.Range("I107").Value = 1248.7
$ws.Range("K107").Value = 1248.7
$ws.Range("M107").Value = 671.3

$ws.Range("H134").Value = 2610.9
$ws.Range("J134").Value = 2949.5
$ws.Range("L134").Value = 8848.5
$ws.Range("N134").Value = -13918.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2059.8
$ws.Range("I16").Value = 1949.75
$ws.Range("J16").Value = 2500
$ws.Range("K16").Value = 1949.75
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = -1662.75
$ws.Range("N16").Value = -3074

$ws.Range("H31").Value = 2009.8
$ws.Range("I31").Value = 1787.5
$ws.Range("J31").Value = 2899
$ws.Range("K31").Value = 1787.5
$ws.Range("L31").Value = 2899
$ws.Range("M31").Value = -1492.5
$ws.Range("N31").Value = -3489

$ws.Range("H34").Value = 2009.8
$ws.Range("I34").Value = 1787.5
$ws.Range("J34").Value = 2899
$ws.Range("K34").Value = 1787.5
$ws.Range("L34").Value = 2899
$ws.Range("M34").Value = -1585.5
$ws.Range("N34").Value = -3303

$ws.Range("H113").Value = 2059.8
$ws.Range("I113").Value = 1949.75
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 1949.75
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 220.25
$ws.Range("N113").Value = -6840

$ws.Range("H122").Value = 4295.778
$ws.Range("I122").Value = 3332.75
$ws.Range("K122").Value = 9998.25
$ws.Range("M122").Value = -7548.25

$ws.Range("H132").Value = 2599.1
$ws.Range("I132").Value = 1832.5
$ws.Range("J132").Value = 3749
$ws.Range("K132").Value = 5497.5
$ws.Range("L132").Value = 11247
$ws.Range("M132").Value = -2967.5
$ws.Range("N132").Value = -16307

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 13025
$ws.Range("I132").Value = 2050
$ws.Range("J132").Value = 24000
$ws.Range("K132").Value = 18450
$ws.Range("L132").Value = 216000
$ws.Range("M132").Value = -15920
$ws.Range("N132").Value = -221060

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 1966.3334
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 1966.3334
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 1966.3334
$ws.Range("N70").Value = -2506.3334
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 1966.3334
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 1966.3334
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 1966.3334
$ws.Range("N73").Value = -3838.3334
$ws.Range("M73").ClearContents()

$ws.Range("H97").Value = 6793.3335
$ws.Range("I97").Value = 6793.3335
$ws.Range("K97").Value = 6793.3335
$ws.Range("M97").Value = -6297.3335

$ws.Range("H122").Value = 2712.8333
$ws.Range("I122").Value = 3407.125
$ws.Range("J122").Value = 1324.25
$ws.Range("K122").Value = 10221.375
$ws.Range("L122").Value = 3972.75
$ws.Range("M122").Value = -7771.375
$ws.Range("N122").Value = -8872.75

$ws.Range("H132").Value = 2153.5625
$ws.Range("I132").Value = 1497.0769
$ws.Range("K132").Value = 4491.2307
$ws.Range("M132").Value = -1961.2307

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1109.5
$ws.Range("I16").Value = 1109.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1109.5
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -939.5
$ws.Range("N16").ClearContents()

$ws.Range("H46").Value = 2720.923
$ws.Range("I46").Value = 1517.1
$ws.Range("J46").Value = 6733.6665
$ws.Range("K46").Value = 1517.1
$ws.Range("L46").Value = 6733.6665
$ws.Range("M46").Value = -1329.1
$ws.Range("N46").Value = -7109.6665

$ws.Range("H55").Value = 1471.2354
$ws.Range("I55").Value = 2367.889
$ws.Range("J55").Value = 462.5
$ws.Range("K55").Value = 2367.889
$ws.Range("L55").Value = 462.5
$ws.Range("M55").Value = -2194.889
$ws.Range("N55").Value = -808.5

$ws.Range("H61").Value = 11899.8
$ws.Range("I61").Value = 9874.75
$ws.Range("K61").Value = 9874.75
$ws.Range("M61").Value = -9672.75

$ws.Range("H86").Value = 10000
$ws.Range("J86").Value = 10000
$ws.Range("L86").Value = 10000
$ws.Range("N86").Value = -12372

$ws.Range("H89").Value = 10000
$ws.Range("J89").Value = 10000
$ws.Range("L89").Value = 30000
$ws.Range("N89").Value = -41856

$ws.Range("H113").Value = 11899.8
$ws.Range("I113").Value = 9874.75
$ws.Range("K113").Value = 9874.75
$ws.Range("M113").Value = -7704.75

$ws.Range("H132").Value = 2627.8215
$ws.Range("I132").Value = 2084.9048
$ws.Range("K132").Value = 6254.714399999999
$ws.Range("M132").Value = -3724.714399999999

$ws.Range("H136").Value = 4138.4287
$ws.Range("I136").Value = 4162.1665
$ws.Range("K136").Value = 12486.4995
$ws.Range("M136").Value = -9936.499500000002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1476.2941
$ws.Range("I100").Value = 1041.4166
$ws.Range("K100").Value = 2082.8332
$ws.Range("M100").Value = -1541.8332

$ws.Range("H132").Value = 2197.2856
$ws.Range("I132").Value = 1097.25
$ws.Range("K132").Value = 3291.75
$ws.Range("M132").Value = -761.75

$ws.Range("H136").Value = 1716.4482
$ws.Range("I136").Value = 1838
$ws.Range("J136").Value = 663
$ws.Range("K136").Value = 5514
$ws.Range("L136").Value = 1989
$ws.Range("M136").Value = -2964
$ws.Range("N136").Value = -7089
